$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Iran -> Iran
$ws.Range("B13").Value = 182525
$ws.Range("C13").Value = 2369
$ws.Range("D13").Value = 144649
$ws.Range("E13").Value = 29217
$ws.Range("G13").Value = 75
$ws.Range("H13").Value = 8659

# Row 23: Catar -> Catar
$ws.Range("B23").Value = 76588
$ws.Range("C23").Value = 1517
$ws.Range("D23").Value = 53296
$ws.Range("E23").Value = 23222
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 70

# Row 24: Belgica -> Belgica
$ws.Range("B24").Value = 59819
$ws.Range("C24").Value = 108
$ws.Range("D24").Value = 16498
$ws.Range("E24").Value = 33675
$ws.Range("G24").Value = 10
$ws.Range("H24").Value = 9646

# Row 34: Portugal -> Indonesia
$ws.Range("A34").Value = "Indonesia"
$ws.Range("B34").Value = 36406
$ws.Range("C34").Value = 1111
$ws.Range("D34").Value = 13213
$ws.Range("E34").Value = 21145
$ws.Range("G34").Value = 48
$ws.Range("H34").Value = 2048

# Row 35: Indonesia -> Portugal
$ws.Range("A35").Value = "Portugal"
$ws.Range("B35").Value = 35910
$ws.Range("D35").Value = 22002
$ws.Range("E35").Value = 12404
$ws.Range("H35").Value = 1504

# Row 36: Kuwait -> Kuwait
$ws.Range("B36").Value = 34952
$ws.Range("C36").Value = 520
$ws.Range("D36").Value = 25048
$ws.Range("E36").Value = 9619
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 285

# Row 37: Suiza -> Suiza
$ws.Range("B37").Value = 31063
$ws.Range("C37").Value = 19
$ws.Range("E37").Value = 325
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 1938

# Row 42: Filipinas -> Filipinas
$ws.Range("B42").Value = 24787
$ws.Range("C42").Value = 612
$ws.Range("D42").Value = 5454
$ws.Range("E42").Value = 18281
$ws.Range("G42").Value = 16
$ws.Range("H42").Value = 1052

# Row 45: Rumania -> Rumania
$ws.Range("B45").Value = 21404
$ws.Range("C45").Value = 222
$ws.Range("D45").Value = 15445
$ws.Range("E45").Value = 4579
$ws.Range("G45").Value = 11
$ws.Range("H45").Value = 1380

# Row 46: Oman -> Oman
$ws.Range("B46").Value = 21071
$ws.Range("C46").Value = 1117
$ws.Range("D46").Value = 7489
$ws.Range("E46").Value = 13486
$ws.Range("G46").Value = 7
$ws.Range("H46").Value = 96

# Row 60: Moldavia -> Moldavia
$ws.Range("D60").Value = 6229
$ws.Range("E60").Value = 4116
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 382

# Row 63: Chequia -> Chequia
$ws.Range("B63").Value = 9887
$ws.Range("C63").Value = 32
$ws.Range("D63").Value = 7170
$ws.Range("E63").Value = 2388
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 329

# Row 67: Guatemala -> Marruecos
$ws.Range("A67").Value = "Marruecos"
$ws.Range("B67").Value = 8581
$ws.Range("C67").Value = 44
$ws.Range("D67").Value = 7600
$ws.Range("E67").Value = 769
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 212

# Row 68: Marruecos -> Guatemala
$ws.Range("A68").Value = "Guatemala"
$ws.Range("B68").Value = 8561
$ws.Range("C68").Value = 340
$ws.Range("D68").Value = 1567
$ws.Range("E68").Value = 6660
$ws.Range("G68").Value = 18
$ws.Range("H68").Value = 334

# Row 69: Malasia -> Malasia
$ws.Range("B69").Value = 8402
$ws.Range("C69").Value = 33
$ws.Range("D69").Value = 7168
$ws.Range("E69").Value = 1115
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 119

# Row 71: Australia -> Australia
$ws.Range("B71").Value = 7290
$ws.Range("C71").Value = 5
$ws.Range("D71").Value = 6783
$ws.Range("E71").Value = 405

# Row 72: Finlandia -> Finlandia
$ws.Range("B72").Value = 7073
$ws.Range("C72").Value = 9
$ws.Range("E72").Value = 548

# Row 74: Tayikistan -> Nepal
$ws.Range("A74").Value = "Nepal"
$ws.Range("B74").Value = 5062
$ws.Range("C74").Value = 448
$ws.Range("D74").Value = 877
$ws.Range("E74").Value = 4169
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 16

# Row 75: Uzbekistan -> Senegal
$ws.Range("A75").Value = "Senegal"
$ws.Range("B75").Value = 4851
$ws.Range("C75").Value = 92
$ws.Range("D75").Value = 3100
$ws.Range("E75").Value = 1695
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 56

# Row 76: Senegal -> Tayikistan
$ws.Range("A76").Value = "Tayikistan"
$ws.Range("B76").Value = 4834
$ws.Range("D76").Value = 3062
$ws.Range("E76").Value = 1723
$ws.Range("H76").Value = 49

# Row 77: Nepal -> Uzbekistan
$ws.Range("A77").Value = "Uzbekistan"
$ws.Range("B77").Value = 4819
$ws.Range("C77").Value = 78
$ws.Range("D77").Value = 3700
$ws.Range("E77").Value = 1100
$ws.Range("H77").Value = 19

# Row 86: Gabon -> El Salvador
$ws.Range("A86").Value = "El Salvador"
$ws.Range("B86").Value = 3481
$ws.Range("C86").Value = 108
$ws.Range("D86").Value = 1587
$ws.Range("E86").Value = 1826
$ws.Range("G86").Value = 4
$ws.Range("H86").Value = 68

# Row 87: El Salvador -> Gabon
$ws.Range("A87").Value = "Gabon"
$ws.Range("B87").Value = 3463
$ws.Range("D87").Value = 1024
$ws.Range("E87").Value = 2416
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 23

# Row 114: Libano -> Albania
$ws.Range("A114").Value = "Albania"
$ws.Range("B114").Value = 1416
$ws.Range("C114").Value = 31
$ws.Range("D114").Value = 1034
$ws.Range("E114").Value = 346
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 36

# Row 115: Guinea-Bisau -> Libano
$ws.Range("A115").Value = "Libano"
$ws.Range("B115").Value = 1402
$ws.Range("D115").Value = 845
$ws.Range("E115").Value = 526
$ws.Range("H115").Value = 31

# Row 116: Albania -> Guinea-Bisau
$ws.Range("A116").Value = "Guinea-Bisau"
$ws.Range("B116").Value = 1389
$ws.Range("D116").Value = 153
$ws.Range("E116").Value = 1224
$ws.Range("H116").Value = 12

# Row 118: Paraguay -> Madagascar
$ws.Range("A118").Value = "Madagascar"
$ws.Range("B118").Value = 1240
$ws.Range("C118").Value = 37
$ws.Range("D118").Value = 344
$ws.Range("E118").Value = 886
$ws.Range("H118").Value = 10

# Row 119: Madagascar -> Paraguay
$ws.Range("A119").Value = "Paraguay"
$ws.Range("B119").Value = 1230
$ws.Range("D119").Value = 624
$ws.Range("E119").Value = 595
$ws.Range("H119").Value = 11

# Row 121: Hong Kong -> Hong Kong
$ws.Range("B121").Value = 1109
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 1060
$ws.Range("E121").Value = 45

# Row 137: Uganda -> Uganda
$ws.Range("B137").Value = 686
$ws.Range("C137").Value = 7
$ws.Range("E137").Value = 525

# Row 139: Malta -> Malta
$ws.Range("B139").Value = 645
$ws.Range("C139").Value = 5
$ws.Range("E139").Value = 36

# Row 165: Gibraltar -> Gibraltar
$ws.Range("D165").Value = 172
$ws.Range("E165").Value = 4
